# NIT-9002698693.xlsx — "Estado de Cuenta" update
# - Adds a new mora period (2507) on top of the existing (now-reversed /
#   descending) period table, so the period list grows from 21 to 22 rows.
# - Refreshes "Valor Mora" (E11) and "Cant. Periodos" (F13) accordingly.
# - Pushes the footer (the underline + "NOMBRE/FIRMA DEL REPRESENTANTE
#   LEGAL" block) down by the one extra row the new period added.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Make room for the extra period row: duplicate the bottom ("closing
#    border") row of the table down one row, then relax the row that used
#    to be the bottom row back to a regular interior row (matching its
#    neighbour above).
# ---------------------------------------------------------------------
$ws.Range("B36:J36").Copy() | Out-Null
$ws.Range("B37:J37").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

$ws.Range("B35:J35").Copy() | Out-Null
$ws.Range("B36:J36").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 2) Shift the footer block (signature line) down one row: row 41 -> 42,
#    row 42 -> 43.
# ---------------------------------------------------------------------
$ws.Range("B41:C41").Copy() | Out-Null
$ws.Range("B42:C42").PasteSpecial(-4122) | Out-Null
$ws.Range("H41:J41").Copy() | Out-Null
$ws.Range("H42:J42").PasteSpecial(-4122) | Out-Null

$ws.Range("B42:C42").Copy() | Out-Null
$ws.Range("B43:C43").PasteSpecial(-4122) | Out-Null
$ws.Range("H42:J42").Copy() | Out-Null
$ws.Range("H43:J43").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("B43").Value = "NOMBRE DEL REPRESENTANTE LEGAL"
$ws.Range("H43").Value = "FIRMA DEL REPRESENTANTE LEGAL"
$ws.Range("B42").Value = "___________________________________"
$ws.Range("H42").Value = "___________________________________"
$ws.Range("B41").Value = ""
$ws.Range("H41").Value = ""

# ---------------------------------------------------------------------
# 3) Rewrite the worker/period table (rows 16-37) newest period first.
# ---------------------------------------------------------------------
$periods = @("2507","2506","2505","2504","2503","2502","2501", `
             "2412","2411","2410","2409","2408","2407","2406","2405", `
             "2404","2403","2402","2401","2312","2311","2310")

for ($i = 0; $i -lt $periods.Length; $i++) {
    $r = 16 + $i
    $ws.Range("B$r").Value = "CC"
    $ws.Range("C$r").Value = "92261505"
    $ws.Range("D$r").Value = "DAIRO DAVID HERRERA MORALES"
    $ws.Range("E$r").Value = $periods[$i]
    $ws.Range("F$r").Value = 46400
    $ws.Range("G$r").Value = 1160000
}

# ---------------------------------------------------------------------
# 4) Update the summary figures: Valor Mora and Cant. Periodos.
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 1020800
$ws.Range("F13").Value = 22

$wb.Worksheets.Item(1).UsedRange | Out-Null
